$d = $word.ActiveDocument

function Replace-InRange($rangeStart, $rangeEnd, $old, $new) {
    $rng = $d.Range($rangeStart, $rangeEnd)
    $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 0, $false, $new, 1) | Out-Null
}

function Replace-InCell($table, $row, $col, $old, $new) {
    $cell = $table.Cell($row, $col)
    Replace-InRange $cell.Range.Start $cell.Range.End $old $new
}

# Update the date line (first paragraph, unique in the document)
$dateRange = $d.Paragraphs.Item(1).Range
Replace-InRange $dateRange.Start $dateRange.End "2025-01-03 Friday" "2025-01-04 Saturday"

# Update each division problem in the table, targeting specific cells (by row/col)
# using wdReplaceOne scoped to the cell range to avoid cross-matching duplicate
# values introduced by earlier replacements (e.g. 585÷6= appears twice across edits).
$t = $d.Tables.Item(1)

Replace-InCell $t 1 1 "987÷9=" "600÷4="
Replace-InCell $t 1 2 "961÷2=" "296÷5="
Replace-InCell $t 1 3 "926÷9=" "843÷9="
Replace-InCell $t 1 4 "526÷2=" "334÷8="
Replace-InCell $t 1 5 "491÷2=" "648÷3="
Replace-InCell $t 5 1 "245÷6=" "423÷5="
Replace-InCell $t 5 2 "761÷9=" "747÷7="
Replace-InCell $t 5 3 "492÷8=" "277÷6="
Replace-InCell $t 5 4 "876÷9=" "331÷9="
Replace-InCell $t 5 5 "996÷7=" "565÷5="
Replace-InCell $t 9 1 "821÷2=" "223÷2="
Replace-InCell $t 9 2 "664÷5=" "597÷5="
Replace-InCell $t 9 3 "495÷8=" "254÷2="
Replace-InCell $t 9 4 "437÷2=" "944÷6="
Replace-InCell $t 9 5 "130÷9=" "585÷6="
Replace-InCell $t 13 1 "694÷2=" "421÷3="
Replace-InCell $t 13 2 "585÷6=" "151÷7="
Replace-InCell $t 13 3 "761÷4=" "492÷4="
Replace-InCell $t 13 4 "888÷8=" "492÷7="
Replace-InCell $t 13 5 "782÷2=" "185÷6="
Replace-InCell $t 17 1 "624÷3=" "969÷5="
Replace-InCell $t 17 2 "238÷2=" "796÷2="
Replace-InCell $t 17 3 "615÷2=" "481÷8="
Replace-InCell $t 17 4 "676÷2=" "332÷8="
Replace-InCell $t 17 5 "603÷9=" "684÷5="
